# Add a default header containing the questionnaire number ("Questionnaire 56")
# to the document's single section, matching the target OOXML:
#   <w:headerReference w:type="default" r:id="rId9"/> on <w:sectPr>
#   word/header1.xml: centered "Header"-styled paragraph, Arial 12pt run
#   "Questionnaire 56"

$d = $word.ActiveDocument

$section = $d.Sections.First
$header = $section.Headers.Item(1)   # wdHeaderFooterPrimary = 1 (the "default" header)

# Set paragraph-level formatting first (style + centering) on the header's
# (currently empty) paragraph.
$header.Range.ParagraphFormat.Style = "Header"
$header.Range.ParagraphFormat.Alignment = 1   # wdAlignParagraphCenter

# Insert the heading text as a new run after the (empty) paragraph mark.
$header.Range.InsertAfter("Questionnaire 56")

# Format only the inserted run (exclude the trailing paragraph mark) so the
# font properties land on the run's <w:rPr>, not the paragraph mark's.
$fullRange = $header.Range
$textRange = $fullRange.Duplicate
$textRange.End = $fullRange.End - 1
$textRange.Font.Name = "Arial"
$textRange.Font.Size = 12
